$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" header columns (A:J) to "_FV2410" and the "_new"
# header columns (L:U) to "_FV2504" -- column K holds the "diff" header
# and stays untouched.
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2410")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2504")
}

# Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a proper Excel Table (ListObject).
$rng = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
